$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new record (test1 / tes1t@gmail.com odoo / mec / ahmed) ---
# Shared-string table insertion order follows the column order the row was
# authored in: nom_etudiant, email, encadrant, titre_sujet.
$ws.Range("A3").Value = "test1"

# B3 gets a mailto hyperlink like B2, with the same "Lien hypertexte" style
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:tes1t@gmail.com odoo", "", "", "tes1t@gmail.com odoo")
$ws.Range("B3").Style = $ws.Range("B2").Style

$ws.Range("D3").Value = "ahmed"
$ws.Range("C3").Value = "mec"

# --- Row 1: extra styled (empty) header-like cells J1:N1, with K1:N1 merged ---

# J1: bold dark label cell with a white right border
$j1 = $ws.Range("J1")
$j1.Font.Name = "Arial"
$j1.Font.Size = 10
$j1.Font.Bold = $true
$j1.Font.Color = 0
$j1.Interior.Color = 16777215
$j1.Borders.Item(10).Weight = -4138
$j1.Borders.Item(10).Color = 16777215
$j1.HorizontalAlignment = -4131
$j1.VerticalAlignment = -4108
$j1.WrapText = $true

# K1: grey value cell with a white left border (start of the merged K1:N1 block)
$k1 = $ws.Range("K1")
$k1.Font.Name = "Arial"
$k1.Font.Size = 10
$k1.Font.Color = 4473924
$k1.Interior.Color = 16777215
$k1.Borders.Item(7).Weight = -4138
$k1.Borders.Item(7).Color = 16777215
$k1.VerticalAlignment = -4108
$k1.WrapText = $true

# L1, M1, N1: same grey font/fill, no border
$ws.Range("L1:N1").Font.Name = "Arial"
$ws.Range("L1:N1").Font.Size = 10
$ws.Range("L1:N1").Font.Color = 4473924
$ws.Range("L1:N1").Interior.Color = 16777215
$ws.Range("L1:N1").VerticalAlignment = -4108
$ws.Range("L1:N1").WrapText = $true

$null = $ws.Range("K1:N1").Merge()

# --- Sheet-level bits ---
$null = $ws.Range("C3").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
